# Insert a new data row above row 126 (shifting existing rows 126..188 down to 127..189)
# and populate it with a new Camote price record for "Vega Modelo de Temuco".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 126 (only within the used columns A:R) - this pushes the
# existing data rows 126..188 down by one, to 127..189, without expanding the sheet
# dimension out to the full row width.
$ws.Range("A126:R126").Insert()

# Copy the formatting of the row that is now 127 (the old row 126) into the new row 126,
# so the new row keeps the same styling (e.g. the date number format in column D).
$ws.Range("A127:R127").Copy()
$ws.Range("A126:R126").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row 126 with the new record's values.
$row = 126
$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"
$ws.Cells.Item($row, 4).Value = 45016
$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = 100114002
$ws.Cells.Item($row, 7).Value = "Camote"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 50
$ws.Cells.Item($row, 11).Value = 25000
$ws.Cells.Item($row, 12).Value = 26000
$ws.Cells.Item($row, 13).Value = 25600
$ws.Cells.Item($row, 14).Value = "`$/malla 20 kilos"
$ws.Cells.Item($row, 15).Value = "Perú"
$ws.Cells.Item($row, 16).Value = 1280
$ws.Cells.Item($row, 17).Value = 20
$ws.Cells.Item($row, 18).Value = "Hortaliza"
